# Helper: build a COM "RGB" value (0x00BBGGRR) from hex R,G,B components
# so srgbClr values in OOXML come out the way we expect.
function RGBVal([int]$r, [int]$g, [int]$b) {
    return $b * 65536 + $g * 256 + $r
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Slide background: solid fill FFDF00 -------------------------------
$s.Background.Fill.Solid()
$s.Background.Fill.ForeColor.RGB = (RGBVal 0xFF 0xDF 0x00)

# --- 2. Title shape ("Delicious Hamburger" -> empty) -----------------------
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Paragraphs(1).Delete()

# --- 3. TextBox 2 shape ------------------------------------------------
$box = $s.Shapes.Item(2)
$tr = $box.TextFrame.TextRange

# First (empty) paragraph becomes the new sub-heading text.
$tr.Paragraphs(1).Text = "Introduction to Fast Food"

# Second paragraph: bump size/bold/color and swap the text, without
# disturbing the pre-existing paragraph-level default run formatting
# object more than necessary.
$para2 = $tr.Paragraphs(2)
$para2.Font.Size = 57.6
$para2.Font.Bold = $true
$para2.Font.Color.RGB = (RGBVal 0xFF 0x00 0x00)
[void]$para2.Replace("A hamburger is a sandwich consisting of a cooked patty of ground meat, usually beef, placed inside a sliced bun.", "Pizza and Hamburgers", 0, 0, 0)

# Re-home the box (moved up) and pin the height back to its original
# value -- editing the text re-triggers the shape's auto-fit layout,
# which would otherwise leave a stale computed height behind.
$box.Top = 72
$box.Height = 72

Write-Host "edit.ps1 applied"
